$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 and row 3 were already empty; the "total" row (previously row 4,
# containing the "مجموع" / sum label plus the G/H/I 00:00:00 totals) needs
# to move up to row 3 - i.e. delete the empty row 3 so row 4 shifts into
# its place, fixing the "extra" gap above the totals row.
$ws.Rows.Item(3).Delete()
